$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8934
$ws.Range("I76").Value = 11053.077
$ws.Range("J76").Value = 4998.5713
$ws.Range("K76").Value = 11053.077
$ws.Range("L76").Value = 4998.5713
$ws.Range("M76").Value = -10738.077
$ws.Range("N76").Value = -5628.5713
$ws.Range("H79").Value = 8934
$ws.Range("I79").Value = 11053.077
$ws.Range("J79").Value = 4998.5713
$ws.Range("K79").Value = 11053.077
$ws.Range("L79").Value = 4998.5713
$ws.Range("M79").Value = -9961.076999999999
$ws.Range("N79").Value = -7182.5713
$ws.Range("H82").Value = 5220
$ws.Range("I82").Value = 5220
$ws.Range("K82").Value = 15660
$ws.Range("M82").Value = -15254
$ws.Range("H85").Value = 5220
$ws.Range("I85").Value = 5220
$ws.Range("K85").Value = 15660
$ws.Range("M85").Value = -14256
$ws.Range("H106").Value = 1677.6666
$ws.Range("I106").Value = 1511.875
$ws.Range("K106").Value = 1511.875
$ws.Range("M106").Value = -880.875
$ws.Range("H107").Value = 977.2381
$ws.Range("I107").Value = 977.2381
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 977.2381
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 942.7619
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2182.8538
$ws.Range("I132").Value = 1713.1111
$ws.Range("J132").Value = 5565
$ws.Range("K132").Value = 5139.3333
$ws.Range("L132").Value = 16695
$ws.Range("M132").Value = -2609.3333
$ws.Range("N132").Value = -21755
$ws.Range("H135").Value = 668.03845
$ws.Range("I135").Value = 416.3684
$ws.Range("J135").Value = 1351.1428
$ws.Range("K135").Value = 3747.3156
$ws.Range("L135").Value = 12160.2852
$ws.Range("M135").Value = -1212.3156
$ws.Range("N135").Value = -17230.2852
$ws.Range("H137").Value = 2763.4348
$ws.Range("J137").Value = 1923.0769
$ws.Range("L137").Value = 5769.2307
$ws.Range("N137").Value = -10869.2307
$ws.Range("H138").Value = 9808620
$ws.Range("I138").Value = 1303.7858
$ws.Range("J138").Value = 16673742
$ws.Range("K138").Value = 3911.3574
$ws.Range("L138").Value = 50021226
$ws.Range("M138").Value = 1228.6426
$ws.Range("N138").Value = -50031506

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2013.3125
$ws.Range("J45").Value = 2692.6667
$ws.Range("L45").Value = 2692.6667
$ws.Range("N45").Value = -3446.6667
$ws.Range("H74").Value = 22753500
$ws.Range("I74").Value = 22753500
$ws.Range("K74").Value = 22753500
$ws.Range("M74").Value = -22752626
$ws.Range("H77").Value = 22753500
$ws.Range("I77").Value = 22753500
$ws.Range("K77").Value = 113767500
$ws.Range("M77").Value = -113763132
$ws.Range("H102").Value = 2405.8
$ws.Range("I102").Value = 1609
$ws.Range("K102").Value = 1609
$ws.Range("M102").Value = 13
$ws.Range("H122").Value = 2971.875
$ws.Range("I122").Value = 2535.7144
$ws.Range("J122").Value = 3453.9473
$ws.Range("K122").Value = 7607.1432
$ws.Range("L122").Value = 10361.8419
$ws.Range("M122").Value = -5157.1432
$ws.Range("N122").Value = -15261.8419
$ws.Range("H132").Value = 77034580
$ws.Range("I132").Value = 4392
$ws.Range("K132").Value = 13176
$ws.Range("M132").Value = -10646

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1815.1666
$ws.Range("I20").Value = 1723
$ws.Range("J20").Value = 1907.3334
$ws.Range("K20").Value = 1723
$ws.Range("L20").Value = 1907.3334
$ws.Range("M20").Value = -1476
$ws.Range("N20").Value = -2401.3334
$ws.Range("H134").Value = 3669.0625
$ws.Range("I134").Value = 3400.4614
$ws.Range("K134").Value = 10201.3842
$ws.Range("M134").Value = -7666.3842
$ws.Range("H140").Value = 119080
$ws.Range("J140").Value = 119080
$ws.Range("L140").Value = 119080
$ws.Range("N140").Value = -129440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6652.8
$ws.Range("J31").Value = 9655.929
$ws.Range("L31").Value = 9655.929
$ws.Range("N31").Value = -10245.929
$ws.Range("H34").Value = 6652.8
$ws.Range("J34").Value = 9655.929
$ws.Range("L34").Value = 9655.929
$ws.Range("N34").Value = -10059.929
$ws.Range("H58").Value = 2144.9524
$ws.Range("I58").Value = 1604.1818
$ws.Range("K58").Value = 1604.1818
$ws.Range("M58").Value = -1401.1818
$ws.Range("H122").Value = 2374.8572
$ws.Range("I122").Value = 1940.6666
$ws.Range("K122").Value = 5821.9998
$ws.Range("M122").Value = -3371.9998
$ws.Range("H132").Value = 6380.5884
$ws.Range("I132").Value = 5533.643
$ws.Range("K132").Value = 16600.929
$ws.Range("M132").Value = -14070.929
$ws.Range("H136").Value = 2144.9524
$ws.Range("I136").Value = 1604.1818
$ws.Range("K136").Value = 4812.5454
$ws.Range("M136").Value = -2262.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 626.8889
$ws.Range("J122").Value = 577.5
$ws.Range("L122").Value = 5197.5
$ws.Range("N122").Value = -10097.5
$ws.Range("H131").Value = 23983.32
$ws.Range("J131").Value = 4527.1562
$ws.Range("L131").Value = 13581.4686
$ws.Range("N131").Value = -23661.4686

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7444.9546
$ws.Range("I70").Value = 6811.875
$ws.Range("J70").Value = 9133.166999999999
$ws.Range("K70").Value = 6811.875
$ws.Range("L70").Value = 9133.166999999999
$ws.Range("M70").Value = -6541.875
$ws.Range("N70").Value = -9673.166999999999
$ws.Range("H73").Value = 7444.9546
$ws.Range("I73").Value = 6811.875
$ws.Range("J73").Value = 9133.166999999999
$ws.Range("K73").Value = 6811.875
$ws.Range("L73").Value = 9133.166999999999
$ws.Range("M73").Value = -5875.875
$ws.Range("N73").Value = -11005.167
$ws.Range("H80").Value = 3789.318
$ws.Range("I80").Value = 3781.6667
$ws.Range("J80").Value = 3805.7144
$ws.Range("K80").Value = 3781.6667
$ws.Range("L80").Value = 3805.7144
$ws.Range("M80").Value = -2783.6667
$ws.Range("N80").Value = -5801.7144
$ws.Range("H83").Value = 3789.318
$ws.Range("I83").Value = 3781.6667
$ws.Range("J83").Value = 3805.7144
$ws.Range("K83").Value = 18908.3335
$ws.Range("L83").Value = 19028.572
$ws.Range("M83").Value = -13916.3335
$ws.Range("N83").Value = -29012.572
$ws.Range("H86").Value = 12250
$ws.Range("J86").Value = 12250
$ws.Range("L86").Value = 12250
$ws.Range("N86").Value = -14622
$ws.Range("H89").Value = 12250
$ws.Range("J89").Value = 12250
$ws.Range("L89").Value = 36750
$ws.Range("N89").Value = -48606
$ws.Range("H126").Value = 6557.4
$ws.Range("J126").Value = 6109.778
$ws.Range("L126").Value = 18329.334
$ws.Range("N126").Value = -23269.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4237.5713
$ws.Range("I40").Value = 4110.6665
$ws.Range("K40").Value = 4110.6665
$ws.Range("M40").Value = -3974.6665
$ws.Range("H55").Value = 684.45
$ws.Range("J55").Value = 969.7
$ws.Range("L55").Value = 969.7
$ws.Range("N55").Value = -1315.7
$ws.Range("H93").Value = 1993.2593
$ws.Range("I93").Value = 1307.9166
$ws.Range("J93").Value = 2541.5334
$ws.Range("K93").Value = 1307.9166
$ws.Range("L93").Value = 2541.5334
$ws.Range("M93").Value = -59.91660000000002
$ws.Range("N93").Value = -5037.5334
$ws.Range("H100").Value = 4348
$ws.Range("I100").Value = 3597
$ws.Range("J100").Value = 5099
$ws.Range("K100").Value = 3597
$ws.Range("L100").Value = 5099
$ws.Range("M100").Value = -3056
$ws.Range("N100").Value = -6181
$ws.Range("H132").Value = 1981.5
$ws.Range("I132").Value = 1981.5
$ws.Range("K132").Value = 5944.5
$ws.Range("M132").Value = -3414.5
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -90060
$ws.Range("H136").Value = 2399.8542
$ws.Range("I136").Value = 1774.4103
$ws.Range("K136").Value = 5323.2309
$ws.Range("M136").Value = -2773.2309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7688.778
$ws.Range("I62").Value = 7424.75
$ws.Range("J62").Value = 7900
$ws.Range("K62").Value = 7424.75
$ws.Range("L62").Value = 7900
$ws.Range("M62").Value = -6800.75
$ws.Range("N62").Value = -9148
$ws.Range("H65").Value = 7688.778
$ws.Range("I65").Value = 7424.75
$ws.Range("J65").Value = 7900
$ws.Range("K65").Value = 37123.75
$ws.Range("L65").Value = 39500
$ws.Range("M65").Value = -34003.75
$ws.Range("N65").Value = -45740
$ws.Range("H122").Value = 2396.65
$ws.Range("I122").Value = 1979.0625
$ws.Range("K122").Value = 5937.1875
$ws.Range("M122").Value = -3487.1875
$ws.Range("H132").Value = 3547.4043
$ws.Range("I132").Value = 3359.6155
$ws.Range("J132").Value = 4462.875
$ws.Range("K132").Value = 10078.8465
$ws.Range("L132").Value = 13388.625
$ws.Range("M132").Value = -7548.8465
$ws.Range("N132").Value = -18448.625
$ws.Range("H136").Value = 1031.279
$ws.Range("I136").Value = 846.8857400000001
$ws.Range("K136").Value = 2540.65722
$ws.Range("M136").Value = 9.342779999999948
